$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 26.127733
$ws.Range("H2").Value = 78.383199
$ws.Range("I2").Value = 0.2666992864894373
$ws.Range("J2").Value = 0.2666992864894374
$ws.Range("M2").Value = 15.24491733333333
$ws.Range("N2").Value = 45.73475199999999
$ws.Range("O2").Value = 0.4831257321597052
$ws.Range("P2").Value = 0.4831257321597052
$ws.Range("Q2").Value = 398.3151296924053
$ws.Range("R2").Value = 3584.836167231648
$ws.Range("S2").Value = 0.1288492880516804
$ws.Range("T2").Value = 0.1288492880516804
$ws.Range("G3").Value = 26.127733
$ws.Range("H3").Value = 78.383199
$ws.Range("I3").Value = 0.2666992864894373
$ws.Range("J3").Value = 0.2666992864894374
$ws.Range("O3").Value = 0.327710667227878
$ws.Range("P3").Value = 0.327710667227878
$ws.Range("Q3").Value = 270.182497493856
$ws.Range("R3").Value = 2431.642477444704
$ws.Range("S3").Value = 0.08740020112465249
$ws.Range("T3").Value = 0.08740020112465251
$ws.Range("G4").Value = 26.127733
$ws.Range("H4").Value = 78.383199
$ws.Range("I4").Value = 0.2666992864894373
$ws.Range("J4").Value = 0.2666992864894374
$ws.Range("M4").Value = 5.969012333333333
$ws.Range("N4").Value = 17.907037
$ws.Range("O4").Value = 0.1891636006124168
$ws.Range("P4").Value = 0.1891636006124168
$ws.Range("Q4").Value = 155.9567605190403
$ws.Range("R4").Value = 1403.610844671363
$ws.Range("S4").Value = 0.05044979731310444
$ws.Range("T4").Value = 0.05044979731310446
$ws.Range("I5").Value = 0.2440410104700376
$ws.Range("J5").Value = 0.2440410104700377
$ws.Range("M5").Value = 15.24491733333333
$ws.Range("N5").Value = 45.73475199999999
$ws.Range("O5").Value = 0.4831257321597052
$ws.Range("P5").Value = 0.4831257321597052
$ws.Range("Q5").Value = 364.475016094535
$ws.Range("R5").Value = 3280.275144850815
$ws.Range("S5").Value = 0.1179024918603312
$ws.Range("T5").Value = 0.1179024918603312
$ws.Range("I6").Value = 0.2440410104700376
$ws.Range("J6").Value = 0.2440410104700377
$ws.Range("O6").Value = 0.327710667227878
$ws.Range("P6").Value = 0.327710667227878
$ws.Range("Q6").Value = 247.228294337152
$ws.Range("S6").Value = 0.07997484237210158
$ws.Range("T6").Value = 0.0799748423721016
$ws.Range("I7").Value = 0.2440410104700376
$ws.Range("J7").Value = 0.2440410104700377
$ws.Range("M7").Value = 5.969012333333333
$ws.Range("N7").Value = 17.907037
$ws.Range("O7").Value = 0.1891636006124168
$ws.Range("P7").Value = 0.1891636006124168
$ws.Range("Q7").Value = 142.7069638156217
$ws.Range("R7").Value = 1284.362674340596
$ws.Range("S7").Value = 0.04616367623760482
$ws.Range("T7").Value = 0.04616367623760483
$ws.Range("G8").Value = 47.93131266666666
$ws.Range("H8").Value = 143.793938
$ws.Range("I8").Value = 0.489259703040525
$ws.Range("J8").Value = 0.4892597030405251
$ws.Range("M8").Value = 15.24491733333333
$ws.Range("N8").Value = 45.73475199999999
$ws.Range("O8").Value = 0.4831257321597052
$ws.Range("P8").Value = 0.4831257321597052
$ws.Range("Q8").Value = 730.7088992814861
$ws.Range("R8").Value = 6576.380093533375
$ws.Range("S8").Value = 0.2363739522476936
$ws.Range("T8").Value = 0.2363739522476936
$ws.Range("G9").Value = 47.93131266666666
$ws.Range("H9").Value = 143.793938
$ws.Range("I9").Value = 0.489259703040525
$ws.Range("J9").Value = 0.4892597030405251
$ws.Range("O9").Value = 0.327710667227878
$ws.Range("P9").Value = 0.327710667227878
$ws.Range("Q9").Value = 495.649651825472
$ws.Range("R9").Value = 4460.846866429248
$ws.Range("S9").Value = 0.1603356237311239
$ws.Range("T9").Value = 0.1603356237311239
$ws.Range("G10").Value = 47.93131266666666
$ws.Range("H10").Value = 143.793938
$ws.Range("I10").Value = 0.489259703040525
$ws.Range("J10").Value = 0.4892597030405251
$ws.Range("M10").Value = 5.969012333333333
$ws.Range("N10").Value = 17.907037
$ws.Range("O10").Value = 0.1891636006124168
$ws.Range("P10").Value = 0.1891636006124168
$ws.Range("Q10").Value = 286.1025964601895
$ws.Range("R10").Value = 2574.923368141706
$ws.Range("S10").Value = 0.09255012706170752
$ws.Range("T10").Value = 0.09255012706170755
